$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "22.446.29"
$ws.Range("E2").Value = "  +0.27%  "
Set-TextValue $ws.Range("D3") "1.573.06"
$ws.Range("E3").Value = "  +0.08%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("E5").Value = "  +0.01%  "
Set-TextValue $ws.Range("D6") "291.80"
$ws.Range("E6").Value = "  +0.18%  "
Set-TextValue $ws.Range("D7") "0.3724"
$ws.Range("E7").Value = "  -0.97%  "
Set-TextValue $ws.Range("D8") "49.92"
$ws.Range("E8").Value = "  -0.23%  "
Set-TextValue $ws.Range("D9") "0.3391"
$ws.Range("E9").Value = "  -0.94%  "
Set-TextValue $ws.Range("D10") "0.07567"
$ws.Range("E10").Value = "  -1.00%  "
$ws.Range("E11").Value = "  -0.61%  "
$ws.Range("E12").Value = "  -0.10%  "
Set-TextValue $ws.Range("D13") "21.27"
$ws.Range("E13").Value = "  +0.48%  "
Set-TextValue $ws.Range("D14") "6.014"
$ws.Range("E14").Value = "  +0.14%  "
Set-TextValue $ws.Range("D15") "6.956"
$ws.Range("E15").Value = "  +0.10%  "
Set-TextValue $ws.Range("D16") "1.575.15"
$ws.Range("E16").Value = "  +0.21%  "
$ws.Range("E17").Value = "  -0.89%  "
Set-TextValue $ws.Range("D18") "90.99"
$ws.Range("E18").Value = "  +1.22%  "
Set-TextValue $ws.Range("D19") "0.06764"
$ws.Range("E19").Value = "  +0.32%  "
$ws.Range("E20").Value = "  -0.01%  "
Set-TextValue $ws.Range("D21") "6.303"
$ws.Range("E21").Value = "  +1.46%  "
Set-TextValue $ws.Range("D22") "16.30"
$ws.Range("E22").Value = "  -2.81%  "
$ws.Range("E23").Value = "  +1.33%  "
Set-TextValue $ws.Range("D24") "22.438.15"
$ws.Range("E24").Value = "  +0.27%  "
Set-TextValue $ws.Range("D25") "2.336"
$ws.Range("E25").Value = "  -2.56%  "
Set-TextValue $ws.Range("D26") "2.688"
$ws.Range("E26").Value = "  +0.58%  "
Set-TextValue $ws.Range("D27") "20.07"
$ws.Range("E27").Value = "  -0.66%  "
Set-TextValue $ws.Range("D28") "148.79"
$ws.Range("E28").Value = "  +1.16%  "
Set-TextValue $ws.Range("D29") "5.033"
$ws.Range("E29").Value = "  -0.05%  "
Set-TextValue $ws.Range("D30") "125.54"
$ws.Range("E30").Value = "  -0.63%  "
Set-TextValue $ws.Range("D31") "1.750.71"
$ws.Range("E31").Value = "  +0.22%  "
Set-TextValue $ws.Range("D32") "1.060"
$ws.Range("E32").Value = "  +7.62%  "
Set-TextValue $ws.Range("D33") "6.160"
$ws.Range("E33").Value = "  +0.24%  "
Set-TextValue $ws.Range("D34") "1.988"
$ws.Range("E34").Value = "  -0.99%  "
Set-TextValue $ws.Range("D35") "9.843"
$ws.Range("E35").Value = "  -1.12%  "
Set-TextValue $ws.Range("D36") "0.08358"
$ws.Range("E36").Value = "  -1.49%  "
Set-TextValue $ws.Range("D37") "0.02492"
$ws.Range("E37").Value = "  -2.14%  "
$ws.Range("E38").Value = "  -2.03%  "
Set-TextValue $ws.Range("D39") "0.2305"
$ws.Range("E39").Value = "  -0.40%  "
Set-TextValue $ws.Range("D40") "0.06515"
$ws.Range("E40").Value = "  -0.83%  "
Set-TextValue $ws.Range("D41") "5.463"
$ws.Range("E41").Value = "  +0.79%  "
$ws.Range("E42").Value = "  -1.08%  "
Set-TextValue $ws.Range("D43") "0.6211"
$ws.Range("E43").Value = "  -2.94%  "
Set-TextValue $ws.Range("D45") "13.98"
$ws.Range("E45").Value = "  -0.06%  "
Set-TextValue $ws.Range("D46") "3.815"
$ws.Range("E46").Value = "  +0.72%  "
Set-TextValue $ws.Range("D47") "0.5803"
$ws.Range("E47").Value = "  -2.80%  "
Set-TextValue $ws.Range("D48") "129.92"
$ws.Range("E48").Value = "  +3.69%  "
Set-TextValue $ws.Range("D49") "2.065"
$ws.Range("E49").Value = "  -1.11%  "
Set-TextValue $ws.Range("D50") "1.221"
$ws.Range("E50").Value = "  -5.77%  "
Set-TextValue $ws.Range("D51") "0.07320"
$ws.Range("E51").Value = "  -0.14%  "
